$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vin 10V")
$ws.Activate()

# Clear the footnote annotation that was tied to the removed Va column
$ws.Range("J24").ClearContents()

# Remove the "Va" measurement column (drifting Keithley 177 readings)
$ws.Columns("E:E").Delete()

# Shift the embedded chart picture left by one column to follow the
# deleted column (picture was anchored starting at column J, now column I)
$shp = $ws.Shapes.Item(1)
$shp.Left = 467.5

[void]$ws.Range("A32").Select()
